# Update NATMI TPM-derived ligand/receptor/edge expression-specificity
# figures on Sheet1 (Fgf2-Sdc4 LR-pair table) to reflect the re-run with
# the new TPM values. Only the numeric value columns (G..J, M..T) change;
# the categorical/id columns (A..F, K, L) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.04961540057666666
$ws.Range("R2").Value = 0.44653860519
$ws.Range("S2").Value = 0.00004443166191349298
$ws.Range("T2").Value = 0.00004443166191349299
$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("Q3").Value = 2.776325851313889
$ws.Range("R3").Value = 24.986932661825
$ws.Range("S3").Value = 0.002486259712781237
$ws.Range("T3").Value = 0.002486259712781238
$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 6.821501368016667
$ws.Range("R4").Value = 61.39351231215
$ws.Range("S4").Value = 0.006108801682610724
$ws.Range("T4").Value = 0.006108801682610725
$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 4.645053410596
$ws.Range("R5").Value = 41.805480695364
$ws.Range("S5").Value = 0.004159745569136445
$ws.Range("T5").Value = 0.004159745569136446
$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("Q6").Value = 259.9229617151634
$ws.Range("S6").Value = 0.2327666213363828
$ws.Range("T6").Value = 0.2327666213363828
$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 638.63715351706
$ws.Range("R7").Value = 5747.73438165354
$ws.Range("S7").Value = 0.5719133527223839
$ws.Range("T7").Value = 0.571913352722384
$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("M8").Value = 0.303146
$ws.Range("N8").Value = 0.909438
$ws.Range("O8").Value = 0.005142855213700541
$ws.Range("P8").Value = 0.005142855213700542
$ws.Range("Q8").Value = 1.048191359854667
$ws.Range("R8").Value = 9.433722238692001
$ws.Range("S8").Value = 0.0009386779826506027
$ws.Range("T8").Value = 0.0009386779826506028
$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("O9").Value = 0.2877784259203595
$ws.Range("P9").Value = 0.2877784259203595
$ws.Range("Q9").Value = 58.65357803554556
$ws.Range("R9").Value = 527.88220231991
$ws.Range("S9").Value = 0.05252554487119539
$ws.Range("T9").Value = 0.05252554487119539
$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 41.67881
$ws.Range("N10").Value = 125.03643
$ws.Range("O10").Value = 0.7070787188659401
$ws.Range("P10").Value = 0.7070787188659401
$ws.Range("Q10").Value = 144.1132936968467
$ws.Range("R10").Value = 1297.01964327162
$ws.Range("S10").Value = 0.1290565644609455
$ws.Range("T10").Value = 0.1290565644609455
